# CRUD_matrix.xlsx update — customer responses resolved uncertain "?" answers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yellow = 65535   # RGB(255,255,0)

# --- Program (row 14): C14 "?CRUD" -> "CRUD" (bold+highlight) ---
$c14 = $ws.Range("C14")
$c14.Value() = "CRUD"
$c14.Font.Bold = $true
$c14.Interior.Color = $yellow
$c14.HorizontalAlignment = -4108

# --- Department (row 13): D13 "?R" -> "R", highlighted ---
$d13 = $ws.Range("D13")
$d13.Value() = "R"
$d13.Interior.Color = $yellow
$d13.HorizontalAlignment = -4108

# --- Program (row 14): D14 "?R" -> "R" (highlight) ---
$d14 = $ws.Range("D14")
$d14.Value() = "R"
$d14.Interior.Color = $yellow
$d14.HorizontalAlignment = -4108

# --- Subject (Lesson, Session) (row 15): D15 "?R" -> "R", highlighted ---
$d15 = $ws.Range("D15")
$d15.Value() = "R"
$d15.Interior.Color = $yellow
$d15.HorizontalAlignment = -4108

# --- SubjectEquipment (row 16): D16 "?R" -> "R", highlighted ---
$d16 = $ws.Range("D16")
$d16.Value() = "R"
$d16.Interior.Color = $yellow
$d16.HorizontalAlignment = -4108

# --- User (row 18): D18 "?" -> "R", highlighted ---
$d18 = $ws.Range("D18")
$d18.Value() = "R"
$d18.Interior.Color = $yellow
$d18.HorizontalAlignment = -4108

# --- SpaceType (row 9): G9 "?" cleared, highlighted (no border) ---
$g9 = $ws.Range("G9")
$g9.ClearContents()
$g9.Interior.Color = $yellow
$g9.HorizontalAlignment = -4108

# --- H9: confirmation text updated with customer response date ---
$h9 = $ws.Range("H9")
$h9.Value() = "Confirm from Customer 2024-02-06"

# --- View: zoom back down to normal 100% ---
$excel.ActiveWindow.Zoom = 100
